$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.990.24"
$ws.Range("E2").Value = "  +0.20%  "
$ws.Range("D3").Value = "1.632.77"
$ws.Range("E3").Value = "  -0.64%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.997"
$ws.Range("E4").Value = "  -0.28%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.99"
$ws.Range("E5").Value = "  -0.72%  "
$ws.Range("E6").Value = "  -0.50%  "
$ws.Range("E7").Value = "  -0.30%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.52"
$ws.Range("E8").Value = "  -0.60%  "
$ws.Range("E9").Value = "  -2.28%  "
$ws.Range("E10").Value = "  -0.41%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0880"
$ws.Range("E11").Value = "  +0.78%  "
$ws.Range("D12").Value = "1.865.19"
$ws.Range("E12").Value = "  -0.57%  "
$ws.Range("D13").Value = "1.630.42"
$ws.Range("E13").Value = "  -0.81%  "
$ws.Range("E14").Value = "  -0.27%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.562"
$ws.Range("E15").Value = "  -1.78%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.60"
$ws.Range("E16").Value = "  -0.25%  "
$ws.Range("D17").Value = "27.982.26"
$ws.Range("E17").Value = "  +0.26%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "232.08"
$ws.Range("E18").Value = "  +0.70%  "
$ws.Range("D19").Value = "0.0₃0727"
$ws.Range("E19").Value = "  +0.16%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.54"
$ws.Range("E20").Value = "  -1.03%  "
$ws.Range("E21").Value = "  -0.39%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.47"
$ws.Range("E22").Value = "  -4.29%  "
$ws.Range("E23").Value = "  -0.84%  "
$ws.Range("E24").Value = "  -3.48%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "154.22"
$ws.Range("E25").Value = "  +1.31%  "
$ws.Range("E26").Value = "  +0.23%  "
$ws.Range("E27").Value = "  -0.72%  "
$ws.Range("E28").Value = "  -0.49%  "
$ws.Range("E29").Value = "  -0.26%  "
$ws.Range("E30").Value = "  -0.50%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0481"
$ws.Range("E31").Value = "  -0.85%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.40"
$ws.Range("E32").Value = "  +1.79%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.11"
$ws.Range("E33").Value = "  +0.39%  "
$ws.Range("D34").Value = "1.407.98"
$ws.Range("E34").Value = "  -1.35%  "
$ws.Range("E35").Value = "  -0.39%  "
$ws.Range("E36").Value = "  +8.44%  "
$ws.Range("E37").Value = "  +0.45%  "
$ws.Range("E39").Value = "  -0.22%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.868"
$ws.Range("E40").Value = "  -2.49%  "
$ws.Range("E41").Value = "  -1.08%  "
$ws.Range("E42").Value = "  -0.31%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "66.96"
$ws.Range("E43").Value = "  -2.50%  "
$ws.Range("E44").Value = "  +0.56%  "
$ws.Range("E45").Value = "  +0.28%  "
$ws.Range("E46").Value = "  -0.52%  "
$ws.Range("D47").Value = "1.775.36"
$ws.Range("E47").Value = "  -0.50%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "88.14"
$ws.Range("E48").Value = "  -1.29%  "
$ws.Range("E49").Value = "  -3.60%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0999"
$ws.Range("E50").Value = "  -1.04%  "
$ws.Range("E51").Value = "  -0.23%  "
